$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Life Steal" entry as rows 26-27, mirroring the layout used by
#     every other two-row entry block (data row + blank continuation row). ---

# Copy the formatting of the previous entry block (row 24/25, the "boss1"-style
# two-row record) down onto the two new rows so borders/alignment/wrap match.
$ws.Range("A24:D24").Copy() | Out-Null
$ws.Range("A26:D26").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("A25:D25").Copy() | Out-Null
$ws.Range("A27:D27").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$excel.CutCopyMode = $false

# Row 26: the actual data for the new entry (ID 18, Life Steal)
$ws.Range("A26").Value = 18
$ws.Range("B26").Value = "-"
$ws.Range("C26").Value = "Life Steal"
$ws.Range("D26").Value = "-"

# Row 27 cells stay blank (matches every other continuation row in the sheet)

# Merge each column across the new two-row block, same pattern as A24:A25 etc.
$ws.Range("A26:A27").Merge() | Out-Null
$ws.Range("B26:B27").Merge() | Out-Null
$ws.Range("C26:C27").Merge() | Out-Null
$ws.Range("D26:D27").Merge() | Out-Null

# --- Scroll the sheet view up a bit (was showing row 13 at top, now row 10) ---
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
